$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number need to be forced to
# Text format first, otherwise Excel auto-converts them to numeric values and
# mangles formatting (e.g. "13.90" -> 13.9).

$ws.Range("D2").Value = "42.893.70"
$ws.Range("E2").Value = "  -5.18%  "

$ws.Range("D3").Value = "2.216.41"
$ws.Range("E3").Value = "  -6.38%  "

$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.19"
$ws.Range("E5").Value = "  +0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.37"
$ws.Range("E6").Value = "  -9.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.579"
$ws.Range("E7").Value = "  -8.73%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.555"
$ws.Range("E9").Value = "  -9.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.58"
$ws.Range("E10").Value = "  -10.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.90"
$ws.Range("E11").Value = "  -3.87%  "

$ws.Range("E12").Value = "  -10.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.60"
$ws.Range("E13").Value = "  -10.30%  "

$ws.Range("E14").Value = "  -1.78%  "

$ws.Range("D15").Value = "2.553.80"
$ws.Range("E15").Value = "  -6.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.859"
$ws.Range("E16").Value = "  -12.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.30"
$ws.Range("E17").Value = "  -7.15%  "

$ws.Range("D18").Value = "2.215.35"
$ws.Range("E18").Value = "  -6.31%  "

$ws.Range("D19").Value = "42.792.71"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.90"
$ws.Range("E20").Value = "  -10.26%  "

$ws.Range("D21").Value = "0.0₃0960"
$ws.Range("E21").Value = "  -9.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.51"
$ws.Range("E22").Value = "  -10.64%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.80"
$ws.Range("E23").Value = "  -11.43%  "

$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.17"
$ws.Range("E24").Value = "  -12.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "235.49"
$ws.Range("E25").Value = "  -11.00%  "

$ws.Range("E26").Value = "  -8.25%  "

$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("E28").Value = "  +0.86%  "

$ws.Range("E29").Value = "  -2.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.92"
$ws.Range("E30").Value = "  -11.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.30"
$ws.Range("E31").Value = "  -15.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.50"
$ws.Range("E32").Value = "  -4.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.17"
$ws.Range("E33").Value = "  -9.83%  "

$ws.Range("E34").Value = "  -9.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "151.80"
$ws.Range("E35").Value = "  -10.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.63"
$ws.Range("E36").Value = "  -8.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.16"
$ws.Range("E37").Value = "  +2.29%  "

$ws.Range("E39").Value = "  -1.78%  "

$ws.Range("E40").Value = "  -6.89%  "

$ws.Range("E41").Value = "  -11.41%  "

$ws.Range("E42").Value = "  -9.55%  "

$ws.Range("E43").Value = "  -9.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.52"
$ws.Range("E44").Value = "  +4.37%  "

$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("D46").Value = "1.727.26"
$ws.Range("E46").Value = "  -7.87%  "

$ws.Range("E47").Value = "  -11.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.01"
$ws.Range("E48").Value = "  -15.61%  "

$ws.Range("E50").Value = "  -4.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.80"
$ws.Range("E51").Value = "  -12.12%  "
